# Update with Correct Forecast output
# Rebuilds the single "Sheet1" tab into a 4-tab forecast workbook:
#   1) "Sales vs PO"     - existing data, reshaped with a new "Order Week" column
#   2) "Weekly Growth"   - the PO quantities that used to live on sheet 1, plus a Growth% column
#   3) "Volume Insights" - summary stats of the PO quantities
#   4) "Prediction Info" - the predicted next week PO quantity
#
# Cell formatting (bold/bordered headers, yyyy-mm-dd date cells) is carried
# over by copying the format from the original sheet's own header/date cells
# (A1 / A2) rather than re-building it property-by-property, so the existing
# style entries in the workbook are reused as-is instead of minting new ones.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheets: rename the original tab and add the three new ones, in order.
# ---------------------------------------------------------------------------
$wsSales = $wb.Worksheets.Item(1)

# Grab the two canonical, already-styled cells from the original sheet
# *before* renaming/touching anything else: A1 = bold bordered header style,
# A2 = yyyy-mm-dd date style. We'll stamp every header / date cell we create
# from these two so the workbook's existing style table is reused verbatim.
$headerStyleSrc = $wsSales.Range("A1")
$dateStyleSrc = $wsSales.Range("A2")

$wsSales.Name = "Sales vs PO"

$wsGrowth = $wb.Worksheets.Add($null, $wsSales)
$wsGrowth.Name = "Weekly Growth"

$wsVolume = $wb.Worksheets.Add($null, $wsGrowth)
$wsVolume.Name = "Volume Insights"

$wsPred = $wb.Worksheets.Add($null, $wsVolume)
$wsPred.Name = "Prediction Info"

function Set-HeaderCell($ws, $row, $col, $text) {
    $dest = $ws.Cells.Item($row, $col)
    $headerStyleSrc.Copy($dest)
    $dest.Value = $text
}

function Set-DateCell($ws, $row, $col, $serial) {
    $dest = $ws.Cells.Item($row, $col)
    $dateStyleSrc.Copy($dest)
    $dest.Value = $serial
}

# ---------------------------------------------------------------------------
# 2. "Sales vs PO" sheet: ds (shifted +6 days), y, Order Week (old ds), PO_Requested_Qty (zeroed)
# ---------------------------------------------------------------------------
Set-HeaderCell $wsSales 1 1 "ds"
Set-HeaderCell $wsSales 1 2 "y"
Set-HeaderCell $wsSales 1 3 "Order Week"
Set-HeaderCell $wsSales 1 4 "PO_Requested_Qty"

$salesDs      = @(45571, 45578, 45585, 45592, 45599, 45606, 45613, 45620, 45627, 45634, 45641, 45648, 45655)
$salesY       = @(0, 0, 0, 7, 12, 12, 23, 17, 35, 19, 28, 34, 30)
$salesOrderWk = @(45565, 45572, 45579, 45586, 45593, 45600, 45607, 45614, 45621, 45628, 45635, 45642, 45649)

for ($i = 0; $i -lt $salesDs.Count; $i++) {
    $row = $i + 2
    Set-DateCell $wsSales $row 1 $salesDs[$i]
    $wsSales.Cells.Item($row, 2).Value = $salesY[$i]
    Set-DateCell $wsSales $row 3 $salesOrderWk[$i]
    $wsSales.Cells.Item($row, 4).Value = 0
}

# ---------------------------------------------------------------------------
# 3. "Weekly Growth" sheet: ds, PO_Requested_Qty, Growth%
# ---------------------------------------------------------------------------
Set-HeaderCell $wsGrowth 1 1 "ds"
Set-HeaderCell $wsGrowth 1 2 "PO_Requested_Qty"
Set-HeaderCell $wsGrowth 1 3 "Growth%"

$growthDs  = @(45572, 45586, 45593)
$growthQty = @(708, 36, 72)
$growthPct = @(0, -94.91525423728814, 100)

for ($i = 0; $i -lt $growthDs.Count; $i++) {
    $row = $i + 2
    Set-DateCell $wsGrowth $row 1 $growthDs[$i]
    $wsGrowth.Cells.Item($row, 2).Value = $growthQty[$i]
    $wsGrowth.Cells.Item($row, 3).Value = $growthPct[$i]
}

# ---------------------------------------------------------------------------
# 4. "Volume Insights" sheet: Total / Average / Max / Min PO quantity
# ---------------------------------------------------------------------------
Set-HeaderCell $wsVolume 1 1 "Total_PO_Quantity"
Set-HeaderCell $wsVolume 1 2 "Average_PO_Quantity"
Set-HeaderCell $wsVolume 1 3 "Max_PO_Quantity"
Set-HeaderCell $wsVolume 1 4 "Min_PO_Quantity"

$wsVolume.Cells.Item(2, 1).Value = 816
$wsVolume.Cells.Item(2, 2).Value = 272
$wsVolume.Cells.Item(2, 3).Value = 708
$wsVolume.Cells.Item(2, 4).Value = 36

# ---------------------------------------------------------------------------
# 5. "Prediction Info" sheet: Predicted_Next_Week_PO_Quantity
# ---------------------------------------------------------------------------
Set-HeaderCell $wsPred 1 1 "Predicted_Next_Week_PO_Quantity"
$wsPred.Cells.Item(2, 1).Value = 0

# ---------------------------------------------------------------------------
# 6. Leave the first sheet active/selected, matching the original workbook view.
# ---------------------------------------------------------------------------
$wsSales.Select()
